# Switch the first contact's first name and gender text (simulating an
# if/else on gender that swaps the displayed name/gender pair).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Jane"
$ws.Range("D2").Value = "Female"
